# "logic interface art wish list"
#
# Restructure the Sheet1 wish-list:
#  - insert a new "Date Completed" column (B); the old notes column
#    shifts from B to C
#  - record a completed date (with date number format) for the first item
#  - add several new wish-list entries, including a new column D for
#    sub-items under "online support for high scores"
#  - bold the two section headers ("THINGS NEEDED" / "NIFTY SHIT TO ADD")
#  - recolor the alternating row fills to theme colors
#  - update the selected cell on Sheet1 and Sheet2
#
# NOTE: shared-string order matters for a faithful rebuild, so new text
# values are written in the same order the original author entered them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------
# 1. Insert a new column B ("Date Completed"); the existing notes
#    column shifts from B to C.
# ---------------------------------------------------------------
$ws1.Columns("B:B").Insert()

# ---------------------------------------------------------------
# 2. New wish-list text (order chosen to reproduce original shared
#    string table ordering).
# ---------------------------------------------------------------
$ws1.Range("C18").Value = "Level select to show picture of map as well as descriptions of map"
$ws1.Range("B1").Value  = "Date Completed"
$ws1.Range("C19").Value = "online support for high scores"
$ws1.Range("D20").Value = "best time"
$ws1.Range("D21").Value = "shortest algorithm"
$ws1.Range("D22").Value = "best combination"

# Completed date for the first item (number format is applied at the
# very end, after fills/fonts, since changing Interior/Font on a date
# cell otherwise resets it back to the auto-assigned short-date format).
$ws1.Range("B2").Value = [datetime]"2009-03-09"

# ---------------------------------------------------------------
# 3. Bold the section header rows.
# ---------------------------------------------------------------
$ws1.Range("B1:C1").Font.Bold = $true
$ws1.Range("C11").Font.Bold = $true

# ---------------------------------------------------------------
# 4. Recolor fills using theme colors:
#    - the "Corwin" assignment row (row 2) takes the theme equivalent
#      of the old red fill
#    - the rest of the list body takes the theme equivalent of the
#      old green fill
# ---------------------------------------------------------------
$ws1.Range("A2:C2").Interior.ThemeColor = 7

$ws1.Range("C3").Interior.ThemeColor = 6
$ws1.Range("A4:C4").Interior.ThemeColor = 6
$ws1.Range("A5:C5").Interior.ThemeColor = 6
$ws1.Range("A6:C6").Interior.ThemeColor = 6
$ws1.Range("A7:C7").Interior.ThemeColor = 6
$ws1.Range("A8:C8").Interior.ThemeColor = 6
$ws1.Range("C12").Interior.ThemeColor = 6
$ws1.Range("C13").Interior.ThemeColor = 6
$ws1.Range("A14:C14").Interior.ThemeColor = 6
$ws1.Range("C15").Interior.ThemeColor = 6
$ws1.Range("C16").Interior.ThemeColor = 6
$ws1.Range("C17").Interior.ThemeColor = 6
$ws1.Range("C18:C19").Interior.ThemeColor = 6
$ws1.Range("D20:D22").Interior.ThemeColor = 6

# ---------------------------------------------------------------
# 5. Size the new "Date Completed" column to fit its contents.
# ---------------------------------------------------------------
$ws1.Columns("B:B").AutoFit()

# ---------------------------------------------------------------
# 6. Selection bookkeeping to match the source edit.
# ---------------------------------------------------------------
$ws1.Range("C19").Select()
$ws2.Range("A1:I12").Select()
